$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.182.01'
$ws.Range('E2').Value = '  +1.34%  '

$ws.Range('D3').Value = '2.475.51'
$ws.Range('E3').Value = '  +2.55%  '

$ws.Range('E4').Value = '  -0.33%  '

$ws.Range('D5').Value = "'577.11"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.83%  '

$ws.Range('D6').Value = "'146.73"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.60%  '

$ws.Range('E7').Value = '  +0.41%  '

$ws.Range('D8').Value = "'0.541"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.63%  '

$ws.Range('D9').Value = '2.475.13'
$ws.Range('E9').Value = '  +1.39%  '

$ws.Range('D10').Value = "'0.111"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.38%  '

$ws.Range('E11').Value = '  +1.46%  '

$ws.Range('D12').Value = "'5.30"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.89%  '

$ws.Range('E13').Value = '  +1.63%  '

$ws.Range('D14').Value = "'29.18"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +8.63%  '

$ws.Range('E15').Value = '  -0.13%  '

$ws.Range('D16').Value = '2.925.59'
$ws.Range('E16').Value = '  +1.93%  '

$ws.Range('D17').Value = '63.154.17'
$ws.Range('E17').Value = '  +1.55%  '

$ws.Range('D18').Value = '2.466.69'
$ws.Range('E18').Value = '  +1.14%  '

$ws.Range('D19').Value = "'8.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.79%  '

$ws.Range('D20').Value = "'11.10"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.08%  '

$ws.Range('D21').Value = "'330.54"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.41%  '

$ws.Range('E22').Value = '  +9.95%  '

$ws.Range('E23').Value = '  +0.06%  '

$ws.Range('D24').Value = "'0.999"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.07%  '

$ws.Range('D25').Value = "'66.53"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.68%  '

$ws.Range('D26').Value = "'669.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.87%  '

$ws.Range('E27').Value = '  +13.77%  '

$ws.Range('D28').Value = '0.0₃0997'
$ws.Range('E28').Value = '  +1.35%  '

$ws.Range('D29').Value = '2.602.59'
$ws.Range('E29').Value = '  +2.80%  '

$ws.Range('D30').Value = "'0.997"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.13%  '

$ws.Range('D31').Value = "'1.47"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.10%  '

$ws.Range('D32').Value = "'8.13"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.46%  '

$ws.Range('D33').Value = "'1.87"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.20%  '

$ws.Range('E34').Value = '  -0.97%  '

$ws.Range('D35').Value = "'1.56"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.14%  '

$ws.Range('D36').Value = "'0.999"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.29%  '

$ws.Range('D37').Value = "'4.80"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.90%  '

$ws.Range('D38').Value = "'5.53"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.49%  '

$ws.Range('D39').Value = "'0.373"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.24%  '

$ws.Range('D40').Value = "'152.86"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.02%  '

$ws.Range('D41').Value = "'18.81"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.27%  '

$ws.Range('D42').Value = "'2.74"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.32%  '

$ws.Range('E43').Value = '  +0.88%  '

$ws.Range('D45').Value = '0.0₆0303'
$ws.Range('E45').Value = '  +7.54%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'149.32"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.00%  '

$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').Value = "'15.16"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +27.52%  '

$ws.Range('D48').Value = "'3.64"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.70%  '

$ws.Range('D49').Value = "'21.02"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.80%  '

$ws.Range('D50').Value = "'0.609"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.24%  '

$ws.Range('E51').Value = '  +0.76%  '
